$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# EMU -> point conversion constant used by the PowerPoint object model
# (1 point = 12700 EMU). Values below are computed from the exact EMU
# offsets/extents in the target OOXML so that AddTextbox's float32
# Left/Top/Width/Height parameters round-trip back to the same EMU.

$labels = @(
    @{ x = 676.8344881889764;  y = 134.18181102362206;  w = 23.7548031496063; h = 29.081259842519685; t = "1" },
    @{ x = 852.0212598425197;  y = 221.01425196850394;  w = 23.7548031496063; h = 29.081259842519685; t = "2" },
    @{ x = 851.9199212598426;  y = 394.41464566929136;  w = 23.7548031496063; h = 29.081259842519685; t = "3" },
    @{ x = 688.7118897637795;  y = 486.67653543307085;  w = 23.7548031496063; h = 29.081259842519685; t = "4" },
    @{ x = 515.5917322834646;  y = 394.41464566929136;  w = 23.7548031496063; h = 29.081259842519685; t = "5" },
    @{ x = 515.5917322834646;  y = 215.90362204724408;  w = 23.7548031496063; h = 29.081259842519685; t = "6" }
)

foreach ($lbl in $labels) {
    $tb = $s.Shapes.AddTextbox(1, $lbl.x, $lbl.y, $lbl.w, $lbl.h)
    $tb.TextFrame.WordWrap = $false
    $tb.TextFrame.AutoSize = 1
    $tb.Fill.Visible = $false
    $tr = $tb.TextFrame.TextRange
    $tr.Text = $lbl.t
    $tr.LanguageID = "en-GB"
}
